# Commit des appréciations M1S1, M1S2, M2S3, M2S4
#
# Adds a new "Appreciations" column (AB) to the sheet, mirroring the
# formatting of the adjacent "Retards" column (AA), fills in the
# per-student appreciation text, and corrects a handful of absence-hour
# totals (Y/Z columns) that were updated at the same time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Create column AB by copying the formatting of column AA (rows 1-32) ---
# This gives AB the same per-row style ids as AA (s3/s8/s12/s15/...) and
# extends the sheet dimension from A1:AA32 to A1:AB32, matching what Excel
# does automatically when you duplicate a column's formatting.
$ws.Range("AA1:AA32").Copy() | Out-Null
$ws.Range("AB1:AB32").PasteSpecial(-4122) | Out-Null

# New column width for AB (~15.11 characters, closest reachable width)
$ws.Columns("AB").ColumnWidth = 14.3

# --- 2. Header for the new column ---
$ws.Range("AB2").Value = "Appreciations"

# --- 3. Corrected absence-hour totals (Y = justified abs., Z = unjustified abs.) ---
$ws.Range("Z6").Value = "30h00"
$ws.Range("Y8").Value = "44h00"
$ws.Range("Z8").Value = "6h00"
$ws.Range("Y9").Value = "16h00"
$ws.Range("Z9").Value = "16h00"
$ws.Range("Z22").Value = "27h30"
$ws.Range("Z23").Value = "52h00"

# --- 4. Per-student semester appreciation text in column AB ---
$ws.Range("AB4").Value = "SEMESTRE NON VALIDÉ :`n   - UE4 - Méthodologie de la recherche`n   - UE4 - Espi Inside`n   - UE SPE - Gestion de Patrimoine"
$ws.Range("AB5").Value = "SEMESTRE NON VALIDÉ :`n   - UE SPE - Evaluation d'actifs Tertiaires et Industriels"
$ws.Range("AB6").Value = "SEMESTRE NON VALIDÉ :`n   - UE4 - Méthodologie de la recherche`n   - UE4 - Espi Inside"
$ws.Range("AB7").Value = "SEMESTRE VALIDÉ"
$ws.Range("AB8").Value = "SEMESTRE VALIDÉ"
$ws.Range("AB9").Value = "SEMESTRE VALIDÉ"
$ws.Range("AB10").Value = "SEMESTRE VALIDÉ"
$ws.Range("AB11").Value = "SEMESTRE VALIDÉ"
$ws.Range("AB12").Value = "SEMESTRE NON VALIDÉ :`n   - UE1 - Economie Immobilière II`n   - UE4 - Méthodologie de la recherche`n   - UE4 - Espi Inside`n   - UE SPE - Due Diligence`n   - UE SPE - Evaluation d'actifs Tertiaires et Industriels`n   - UE SPE - Gestion de Patrimoine"
$ws.Range("AB13").Value = "SEMESTRE NON VALIDÉ :`n   - UE4 - Espi Inside"
$ws.Range("AB14").Value = "SEMESTRE VALIDÉ"
$ws.Range("AB15").Value = "SEMESTRE VALIDÉ"
$ws.Range("AB16").Value = "SEMESTRE NON VALIDÉ :`n   - UE4 - Espi Inside"
$ws.Range("AB17").Value = "SEMESTRE VALIDÉ"
$ws.Range("AB18").Value = "SEMESTRE NON VALIDÉ :`n   - UE SPE - Droit des sûretés et de la Transmission`n   - UE SPE - Due Diligence"
$ws.Range("AB19").Value = "SEMESTRE VALIDÉ"
$ws.Range("AB20").Value = "SEMESTRE NON VALIDÉ :`n   - UE4 - Espi Inside"
$ws.Range("AB21").Value = "SEMESTRE VALIDÉ"
$ws.Range("AB22").Value = "SEMESTRE NON VALIDÉ :`n   - UE4 - Espi Inside`n   - UE SPE - Droit des sûretés et de la Transmission"
$ws.Range("AB23").Value = "SEMESTRE NON VALIDÉ :`n   - UE4 - Espi Inside"
$ws.Range("AB24").Value = "SEMESTRE NON VALIDÉ :`n   - UE SPE - Droit des sûretés et de la Transmission`n   - UE SPE - Due Diligence"
$ws.Range("AB25").Value = "SEMESTRE NON VALIDÉ :`n   - UE4 - Espi Inside"
$ws.Range("AB26").Value = "SEMESTRE NON VALIDÉ :`n   - UE4 - Espi Inside`n   - UE SPE - Gestion de Patrimoine"
$ws.Range("AB27").Value = "SEMESTRE NON VALIDÉ :`n   - UE4 - Espi Inside"
$ws.Range("AB28").Value = "SEMESTRE VALIDÉ"
$ws.Range("AB29").Value = "SEMESTRE VALIDÉ"
$ws.Range("AB30").Value = "SEMESTRE NON VALIDÉ :`n   - UE4 - Espi Inside"
$ws.Range("AB31").Value = "SEMESTRE VALIDÉ"

# Re-fit row heights: assigning the multi-line appreciation text makes the
# engine mark some rows with an explicit (wrapped) height; AutoFit puts the
# affected rows back to the sheet's normal (non-custom) height, matching
# the source file where no row height changed because of this edit.
$ws.Rows("4:32").AutoFit() | Out-Null

# --- 5. Sheet view: scroll one column to the right and move selection ---
$ws.Range("O1").Select() | Out-Null
$ws.Range("AC5").Select() | Out-Null
